# "versión estable de los dos excel"
#
# On the "Formato" sheet only the selected cell changed (cosmetic).
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Formato")
$ws1.Range("L5").Select()

# On the "Registro" sheet several client/driver names were edited and a
# new row of data (row 10) was added. The order of these edits matters
# because it determines the order new entries are appended to the
# workbook's shared string table.
$ws2 = $wb.Worksheets.Item("Registro")

$ws2.Range("I6").Value = "Centro "           # was "centro "
$ws2.Range("C7").Value = "Michelle P Simón"  # was "Herrería San Simón"
$ws2.Range("C6").Value = "Lucy Martinez"     # was "Jorge Gamez"
$ws2.Range("C9").Value = "chamacon"          # was "Mto Camarillo"
$ws2.Range("C8").Value = "6 2 1"             # was "Antonio Galindo"

# New row of data appended below the existing records.
$ws2.Range("B10").Value = 12121
$ws2.Range("C10").Value = "sdsd"
$ws2.Range("I10").Value = "PUEBLANO"

# Leave the active/selected cell on "Registro" as in the final workbook.
$ws2.Range("H7").Select()
